$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = "2019-12-31 00:00:00"

$ws.Range("O2").Value = 46974411.46
$ws.Range("P2").Value = 260.4995780481
$ws.Range("Q2").Value = 166380969.89
$ws.Range("R2").Value = 922.6762210417
$ws.Range("S2").Value = 21365392.62
$ws.Range("T2").Value = 118.483139849
$ws.Range("U2").Value = -13013909.13
$ws.Range("V2").Value = -72.16945847220001
$ws.Range("Y2").Value = 6575886.71
$ws.Range("Z2").Value = 36.4669968182
$ws.Range("AA2").Value = -15922900
$ws.Range("AB2").Value = -88.3014518411
$ws.Range("AC2").Value = 18032432.84
$ws.Range("AD2").Value = 554.7413872757
